# CoMPAS Initial Roadmap - content update
# - Insert "OSIsoft" into the partner list of every slide footer.
# - Rewrite the "customization/automatization" sentence on the
#   "Main principles for the overall architecture" slide.
# - Rename / re-center / reposition the "User web browser HMI" label.
# - Merge the "Specific microservices" caption into a single shorter
#   paragraph and shrink its text box accordingly.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Footer credits: insert "OSIsoft, " right before "RTE" wherever
#    the "..., National Grid, RTE, Schneider Electric, ..." credit
#    line shows up (every content slide's footer placeholder).
# ---------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*National Grid, RTE*") {
                $hit = $tr.Find("RTE")
                [void]$hit.InsertBefore("OSIsoft, ")
            }
        }
    }
}

# ---------------------------------------------------------------
# 2) "Main principles for the overall architecture" slide: reword
#    "... according to specificities of the users" to
#    "... according to the specificities of users".
# ---------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*specificities of the users*") {
                $hit = $tr.Find("specificities of the users")
                $hit.Text = "the specificities "
                $anchor = $tr.Find("the specificities ")
                [void]$anchor.InsertAfter("of users")
            }
        }
    }
}

# ---------------------------------------------------------------
# 3) "User web browser HMI" label -> "User interface", centered and
#    shifted left on the architecture diagram slide.
# ---------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "User web browser HMI") {
                $shape.TextFrame.TextRange.Text = "User interface"
                $shape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
                $shape.Left = 186.34875
                $shape.Top = 89.29314960629921
            }
        }
    }
}

# ---------------------------------------------------------------
# 4) "Specific microservices or tools through API (e.g. multiple
#    vendor IED Configuration tools)" caption -> single shorter
#    paragraph, with the text box shrunk to fit.
# ---------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "Specific microservices or tools through API*") {
                $shape.TextFrame.TextRange.Text = "Specific microservices or API (e.g. vendor-specific IED Configuration tools)"
                $shape.Height = 39.9867
            }
        }
    }
}
